$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.806.01"
$ws.Range("E2").Value = "  -0.07%  "

# Row 3
$ws.Range("D3").Value = "2.080.95"
$ws.Range("E3").Value = "  -0.13%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.48"
$ws.Range("E5").Value = "  -0.34%  "

# Row 6
$ws.Range("E6").Value = "  -0.14%  "

# Row 7
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.39"
$ws.Range("E7").Value = "  -0.97%  "

# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("E9").Value = "  +0.78%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0785"
$ws.Range("E10").Value = "  -0.69%  "

# Row 11
$ws.Range("E11").Value = "  +3.26%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.09"
$ws.Range("E12").Value = "  +2.16%  "

# Row 13
$ws.Range("D13").Value = "2.387.11"
$ws.Range("E13").Value = "  -0.14%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.30"
$ws.Range("E14").Value = "  +0.27%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.780"
$ws.Range("E15").Value = "  +1.67%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.38"
$ws.Range("E16").Value = "  +1.34%  "

# Row 17
$ws.Range("D17").Value = "2.080.93"
$ws.Range("E17").Value = "  -0.22%  "

# Row 18
$ws.Range("D18").Value = "37.741.52"
$ws.Range("E18").Value = "  +0.06%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.14"
$ws.Range("E19").Value = "  -0.84%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.24"
$ws.Range("E20").Value = "  -0.17%  "

# Row 21
$ws.Range("E21").Value = "  +0.41%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.01"
$ws.Range("E22").Value = "  +0.44%  "

# Row 23
$ws.Range("E23").Value = "  -0.12%  "

# Row 24
$ws.Range("E24").Value = "  -0.72%  "

# Row 25
$ws.Range("E25").Value = "  +0.91%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.83"
$ws.Range("E26").Value = "  +8.95%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.12"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.136"
$ws.Range("E28").Value = "  -2.28%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.51"
$ws.Range("E29").Value = "  -0.14%  "

# Row 30
$ws.Range("E30").Value = "  -0.93%  "

# Row 31
$ws.Range("E31").Value = "  +1.05%  "

# Row 32
$ws.Range("E32").Value = "  +0.54%  "

# Row 33
$ws.Range("E33").Value = "  +0.39%  "

# Row 34
$ws.Range("E34").Value = "  -0.76%  "

# Row 35
$ws.Range("E35").Value = "  -1.55%  "

# Row 36
$ws.Range("E36").Value = "  -0.39%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.40"
$ws.Range("E37").Value = "  -2.11%  "

# Row 38
$ws.Range("E38").Value = "  +0.15%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.45"
$ws.Range("E39").Value = "  +0.65%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0236"
$ws.Range("E40").Value = "  +9.57%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.89"
$ws.Range("E41").Value = "  +4.11%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0973"
$ws.Range("E42").Value = "  -1.62%  "

# Row 43
$ws.Range("E43").Value = "  -1.04%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.85"
$ws.Range("E44").Value = "  +5.18%  "

# Row 45
$ws.Range("D45").Value = "1.455.53"
$ws.Range("E45").Value = "  -0.73%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.16"
$ws.Range("E46").Value = "  -1.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.07"
$ws.Range("E47").Value = "  -0.47%  "

# Row 48
$ws.Range("E48").Value = "  -7.64%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.34"
$ws.Range("E49").Value = "  -0.99%  "

# Row 50
$ws.Range("E50").Value = "  -1.42%  "

# Row 51
$ws.Range("D51").Value = "2.272.01"
$ws.Range("E51").Value = "  -0.14%  "
